$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reg_obl_city")

for ($r = 2; $r -le 127; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}

$ws.Range("A2:A127").Select()
$excel.ActiveWindow.ScrollRow = 106
